# Remove the T_evap column (data column C) from the activity data table.
# This shifts every column from D onward one position to the left
# (headers, units, data values and the eta_oi / Q_amb formulas all move
# with it, and Excel automatically re-points the formulas/styles).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("C").Delete()

# Restore the cell selection left behind by the editing session.
$null = $ws.Range("G31").Select()
